# Add the new veg types (eucalypt, lawn_10cm, lawn_5cm, lawn_2cm) to the
# "updated" sheet of the allometric workbook, as rows 6-9.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("updated")

# Row 6: eucalypt
$ws.Range("A6").Value = "eucalypt"
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = 8.74
$ws.Range("D6").Value = 1.32
$ws.Range("E6").Value = 1.3
$ws.Range("F6").Value = 0.3
$ws.Range("G6").Value = 0.73
$ws.Range("H6").Value = 28.7
$ws.Range("I6").Value = 139.7
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 333.33
$ws.Range("B6").NumberFormat = "0"

# Row 7: lawn_10cm
$ws.Range("A7").Value = "lawn_10cm"
$ws.Range("B7").Value = 31
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = 1.2
$ws.Range("E7").Value = 2.2
$ws.Range("F7").Value = 0.22
$ws.Range("G7").Value = 0.16
$ws.Range("H7").Value = 27.7
$ws.Range("I7").Value = 47.8
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 333.33

# Row 8: lawn_5cm
$ws.Range("A8").Value = "lawn_5cm"
$ws.Range("B8").Value = 32
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = 1.2
$ws.Range("E8").Value = 2.2
$ws.Range("F8").Value = 0.22
$ws.Range("G8").Value = 0.16
$ws.Range("H8").Value = 27.7
$ws.Range("I8").Value = 47.8
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 333.33

# Row 9: lawn_2cm
$ws.Range("A9").Value = "lawn_2cm"
$ws.Range("B9").Value = 33
$ws.Range("C9").Value = 17
$ws.Range("D9").Value = 1.2
$ws.Range("E9").Value = 2.2
$ws.Range("F9").Value = 0.22
$ws.Range("G9").Value = 0.16
$ws.Range("H9").Value = 27.7
$ws.Range("I9").Value = 47.8
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 333.33

# Leave the selection where a user would land after typing the last row
# (one row below the newly-entered data).
$ws.Range("B10").Select()
